$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '28.060.22'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").Value = "'" + '1.833.81'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'" + '324.12'
$ws.Range("E5").Value = '  -3.15%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = "'" + '0.4640'
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("D8").Value = "'" + '0.3880'
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("D9").Value = "'" + '0.07855'
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("D10").Value = "'" + '0.9623'
$ws.Range("E10").Value = '  -2.25%  '
$ws.Range("D11").Value = "'" + '21.92'
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").Value = "'" + '1.821.99'
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("D13").Value = "'" + '5.688'
$ws.Range("E13").Value = '  -2.79%  '
$ws.Range("D14").Value = "'" + '6.917'
$ws.Range("E14").Value = '  -1.36%  '
$ws.Range("D15").Value = "'" + '0.06837'
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = "'" + '88.35'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = "'" + '1.001'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = "'" + '0.000009951'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = "'" + '16.71'
$ws.Range("E19").Value = '  -2.58%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = "'" + '28.072.18'
$ws.Range("E21").Value = '  -2.03%  '
$ws.Range("D22").Value = "'" + '5.312'
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("D23").Value = "'" + '11.01'
$ws.Range("E23").Value = '  -2.92%  '
$ws.Range("D24").Value = "'" + '2.096'
$ws.Range("E24").Value = '  -1.88%  '
$ws.Range("D25").Value = "'" + '2.068.06'
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("D26").Value = "'" + '154.88'
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("D27").Value = "'" + '19.17'
$ws.Range("E27").Value = '  -1.49%  '
$ws.Range("D28").Value = "'" + '5.677'
$ws.Range("E28").Value = '  -5.81%  '
$ws.Range("D29").Value = "'" + '1.961'
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("D30").Value = "'" + '118.24'
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("D31").Value = "'" + '0.9362'
$ws.Range("E31").Value = '  -4.23%  '
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("D33").Value = "'" + '5.259'
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("D34").Value = "'" + '1.321'
$ws.Range("E34").Value = '  -2.23%  '
$ws.Range("D35").Value = "'" + '3.308'
$ws.Range("E35").Value = '  -5.02%  '
$ws.Range("D36").Value = "'" + '0.05872'
$ws.Range("E36").Value = '  -4.62%  '
$ws.Range("D37").Value = "'" + '0.02128'
$ws.Range("E37").Value = '  -3.22%  '
$ws.Range("D38").Value = "'" + '1.145'
$ws.Range("E38").Value = '  -2.16%  '
$ws.Range("D39").Value = "'" + '7.766'
$ws.Range("E39").Value = '  +1.87%  '
$ws.Range("D40").Value = "'" + '0.5603'
$ws.Range("E40").Value = '  -2.15%  '
$ws.Range("D41").Value = "'" + '9.896'
$ws.Range("E41").Value = '  -3.11%  '
$ws.Range("D42").Value = "'" + '0.1764'
$ws.Range("E42").Value = '  -2.16%  '
$ws.Range("D43").Value = "'" + '0.07281'
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("D44").Value = "'" + '11.71'
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("D45").Value = "'" + '0.5277'
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("D46").Value = "'" + '1.155'
$ws.Range("E46").Value = '  -7.29%  '
$ws.Range("D47").Value = "'" + '2.113'
$ws.Range("E47").Value = '  -11.00%  '
$ws.Range("D48").Value = "'" + '1.826'
$ws.Range("E48").Value = '  -4.37%  '
$ws.Range("D49").Value = "'" + '112.20'
$ws.Range("E49").Value = '  -3.06%  '
$ws.Range("D50").Value = "'" + '1.000'
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").Value = "'" + '1.026'
$ws.Range("E51").Value = '  +0.41%  '
